$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1395025849342346, 0.9664603471755981, 0.01732644811272621, 0.9978384375572205),
    @(0.02597234398126602, 0.9944161772727966, 0.01285681594163179, 0.9981797337532043),
    @(0.01831133104860783, 0.9941774010658264, 0.00878625363111496, 0.9985210299491882),
    @(0.01359660737216473, 0.9949120879173279, 0.005402940325438976, 0.9987486004829407),
    @(0.01151106134057045, 0.995316207408905, 0.004714337643235922, 0.9993174076080322),
    @(0.01155757997184992, 0.9953896999359131, 0.004734234418720007, 0.9994311928749084),
    @(0.01037962362170219, 0.9954080581665039, 0.005438288208097219, 0.9994311928749084),
    @(0.01127995643764734, 0.9948018789291382, 0.005241305101662874, 0.9994311928749084),
    @(0.01073058973997831, 0.9949672222137451, 0.004750982392579317, 0.9995449185371399),
    @(0.01088981982320547, 0.9953345656394958, 0.005035346373915672, 0.9995449185371399),
    @(0.009973933920264244, 0.9952243566513062, 0.005365626886487007, 0.9995449185371399),
    @(0.01005701813846827, 0.9954999089241028, 0.005820718593895435, 0.9994311928749084),
    @(0.009635179303586483, 0.9962896704673767, 0.005509081296622753, 0.9995449185371399),
    @(0.008531605824828148, 0.9976673126220703, 0.006082512438297272, 0.9987486004829407),
    @(0.007129268255084753, 0.9978693127632141, 0.00595078943297267, 0.9995449185371399),
    @(0.00708261551335454, 0.9981265068054199, 0.006320483051240444, 0.9994311928749084),
    @(0.007292637601494789, 0.9977958798408508, 0.00592965679243207, 0.9994311928749084),
    @(0.006568565033376217, 0.9981080889701843, 0.006541531998664141, 0.9994311928749084),
    @(0.007550589274615049, 0.9977223873138428, 0.006911433301866055, 0.9993174076080322),
    @(0.006698807701468468, 0.9981265068054199, 0.006295342929661274, 0.9994311928749084),
    @(0.007235968019813299, 0.9977591037750244, 0.006155957467854023, 0.9995449185371399),
    @(0.0072016641497612, 0.9979428052902222, 0.00726740388199687, 0.9993174076080322),
    @(0.007378552574664354, 0.9979428052902222, 0.006112142466008663, 0.9994311928749084),
    @(0.006297868210822344, 0.9980713725090027, 0.006399996113032103, 0.9995449185371399),
    @(0.007147795986384153, 0.9978876709938049, 0.004925710614770651, 0.9994311928749084),
    @(0.00653363298624754, 0.9980162978172302, 0.004610556177794933, 0.9995449185371399),
    @(0.005834972485899925, 0.9984387159347534, 0.005273810122162104, 0.9994311928749084),
    @(0.006740238983184099, 0.998034656047821, 0.004825949668884277, 0.9995449185371399),
    @(0.006428609136492014, 0.9981815814971924, 0.005068634636700153, 0.9995449185371399),
    @(0.006430354434996843, 0.9979428052902222, 0.005155091639608145, 0.9995449185371399),
    @(0.006752344314008951, 0.9979978799819946, 0.004562243353575468, 0.9995449185371399),
    @(0.007262278813868761, 0.9977958798408508, 0.005167745985090733, 0.9995449185371399),
    @(0.006900892127305269, 0.9980530142784119, 0.004928939044475555, 0.9995449185371399),
    @(0.007021408528089523, 0.9978876709938049, 0.005121841095387936, 0.9995449185371399),
    @(0.006412471644580364, 0.9979060888290405, 0.005201456137001514, 0.9997724890708923),
    @(0.006987139116972685, 0.997704029083252, 0.00593773927539587, 0.9995449185371399),
    @(0.006522070616483688, 0.9981448650360107, 0.005954241845756769, 0.9996587038040161),
    @(0.007356030400842428, 0.9977591037750244, 0.006488344166427851, 0.9997724890708923),
    @(0.006055481266230345, 0.9981448650360107, 0.006335740443319082, 0.9997724890708923),
    @(0.00661418167874217, 0.9980713725090027, 0.005276686977595091, 0.9997724890708923),
    @(0.006881479639559984, 0.998034656047821, 0.004894997458904982, 0.9995449185371399),
    @(0.006541045848280191, 0.9981999397277832, 0.004529612604528666, 0.9995449185371399),
    @(0.007214405573904514, 0.9977958798408508, 0.004190125036984682, 0.9997724890708923),
    @(0.007527525536715984, 0.9977591037750244, 0.005610976368188858, 0.9997724890708923),
    @(0.006650392431765795, 0.9978876709938049, 0.005046019796282053, 0.9997724890708923),
    @(0.006275233346968889, 0.9980530142784119, 0.005442051682621241, 0.9997724890708923),
    @(0.006317498162388802, 0.997961163520813, 0.006495713256299496, 0.9995449185371399),
    @(0.006438801996409893, 0.9979795217514038, 0.006329809781163931, 0.9997724890708923),
    @(0.007761974353343248, 0.9979244470596313, 0.005370510276407003, 0.9995449185371399),
    @(0.006703986320644617, 0.9978693127632141, 0.00586738670244813, 0.9997724890708923)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($row, $c).Value = $rowValues[$c - 1]
    }
}
